# Insert a new data row at row 76 (pushing existing rows 76-93 down to 77-94)
# and populate it with the new "Dina" variety record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(76).Insert()

$ws.Range("A76").Value2 = 10
$ws.Range("B76").Value2 = "Vega Modelo de Temuco"
$ws.Range("C76").Value2 = "La Araucanía"
$ws.Range("D76").Value2 = 45258
$ws.Range("E76").Value2 = 9
$ws.Range("F76").Value2 = "Fruta"
$ws.Range("G76").Value2 = 100103
$ws.Range("H76").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I76").Value2 = 100103003
$ws.Range("J76").Value2 = "Damasco"
$ws.Range("K76").Value2 = "Dina"
$ws.Range("L76").Value2 = "Primera"
$ws.Range("M76").Value2 = 55
$ws.Range("N76").Value2 = 36000
$ws.Range("O76").Value2 = 36000
$ws.Range("P76").Value2 = 36000
$ws.Range("Q76").Value2 = "$/caja 15 kilos"
$ws.Range("R76").Value2 = "Región de O'Higgins"
$ws.Range("S76").Value2 = 2400
$ws.Range("T76").Value2 = 15
